# Add a new worksheet "Analysing Business and Consumer" as the last tab,
# populate it with the question bank, format the header row to match the
# other sheets, and make it the active sheet (matching the authored commit).

$wb = $excel.ActiveWorkbook

$data = @(
  ,@("Question", "A", "B", "C", "D", "Correct Option", "Explanation")
  ,@("Which of the following is not a characteristic of B2B markets?", "Few and large buyers", "Many buying influences", "Price-driven impulse purchases", "Geographically concentrated buyers", "C", "B2B purchases are typically planned and complex, not impulse-driven.")
  ,@("In a buying center, the person who actually signs the purchase order is called the:", "Influencer", "Gatekeeper", "Buyer", "Decider", "C", "The Buyer executes the transaction; the Decider chooses the supplier.")
  ,@("Key Account Management focuses primarily on:", "One-time large orders", "Transactional selling", "Long-term relationship building with high-value customers", "Exclusive use of intermediaries", "C", "It aims for sustainable, high-CLV relationships.")
  ,@("Organizational buying refers to:", "Households buying consumer goods", "Firms establishing needs, evaluating, and selecting suppliers", "Impulse buying in retail settings", "Crowd-funded project contributions", "B", "It’s the formal process firms use to purchase products/services.")
  ,@("One lesson from the Berlin Airport (BER) case is the need for:", "More competitive bidding rounds", "Integrated project management and clear ownership", "Decentralized communication channels", "Minimizing stakeholder involvement", "B", "Misaligned incentives and unclear roles derailed the project.")
  ,@("In the S–O–R model of consumer behavior, “O” stands for:", "Offer", "Organism (consumer)", "Opportunity", "Outcome", "B", "It refers to the consumer’s internal processes.")
  ,@("Which level is not one of the three consumer characteristic levels?", "Cultural", "Social", "Personal", "Political", "D", "The three are cultural, social and personal.")
  ,@("Sensation differs from perception in that sensation is:", "Interpretation of stimuli", "Immediate receptor response to stimuli", "Memory integration", "Post-purchase evaluation", "B", "Sensation is the raw input; perception organizes and interprets it.")
  ,@("An example of multimodal advertising would combine:", "Print ads only", "Audio, visual, and haptic cues", "Single-sense billboards", "Text messages alone", "B", "Multimodal engages multiple senses simultaneously.")
  ,@("According to behavioral learning theory, in classical conditioning a(n):", "Reward follows a desired behavior", "Unconditioned stimulus is paired with a neutral stimulus", "Customer surveys are administered", "Price discounts condition repeat purchase", "B", "A neutral stimulus becomes conditioned via pairing.")
  ,@("Maslow’s hierarchy of needs places “self-actualization” at the:", "Base", "Second level", "Middle", "Top", "D", "It’s the highest level of motivational needs.")
  ,@("Which type of advertising appeals directly to a physiological need?", "Luxury car ads", "Food product ads", "Charity donation appeals", "Investment service ads", "B", "Hunger is a physiological need.")
  ,@("Involvement in the buying process is defined as:", "Frequency of purchase", "Perceived relevance based on needs, values, and interests", "Social proof usage", "Visual appeal of packaging", "B", "Involvement rises with personal relevance.")
  ,@("Which step in the consumer decision-making process comes immediately after problem recognition?", "Purchase decision", "Information search", "Post-purchase behavior", "Evaluation of alternatives", "B", "Shoppers seek data once they see a need.")
  ,@("Consumers may use heuristics during evaluation to:", "Analyze every attribute equally", "Simplify complex choices with mental shortcuts", "Ignore brand entirely", "Always select the cheapest option", "B", "Efficiency under time pressure drives heuristic use.")
  ,@("The lexicographic heuristic chooses the option with:", "Highest score on the single most important attribute", "Average across all attributes", "Worst score eliminated first", "Random selection", "A", "Shoppers rank alternatives by key attribute.")
  ,@("Beliefs differ from attitudes in that beliefs are:", "Emotional evaluations", "Thoughts about product attributes", "Behavioral intentions", "Impulse reactions", "B", "Attitudes are favorable/unfavorable evaluations.")
  ,@("Which is not a type of perceived risk in purchase decisions?", "Functional", "Aesthetic", "Financial", "Time-related", "B", "Aesthetic risk falls under psychological or functional.")
  ,@("Value co-creation post-purchase includes all except:", "Supporting “in use” interactions", "Encouraging network relationships", "Complete elimination of customer support", "Engaging employees and partners", "C", "Support is vital, not eliminated.")
  ,@("Cognitive dissonance after buying can be reduced by:", "Ignoring customer feedback", "Offering money-back guarantees", "Removing all warranties", "Raising price post-purchase", "B", "Guarantees reassure the buyer.")
  ,@("Sustainable consumption in postpurchase behavior emphasizes:", "Single-use products", "Long-term environmental and social value", "Planned obsolescence", "Exclusive cost focus", "B", "Sustainability integrates eco and social considerations.")
  ,@("Habitual and impulse buying are contrasted based on:", "Price only", "Involvement level", "Distribution channel", "Brand equity", "B", "Low involvement yields habitual/impulse choices.")
  ,@("Which source is not among the four information search types?", "Personal", "Commercial", "Public", "Psychological", "D", "The four are personal, commercial, public, experiential.")
  ,@("Consumers select more familiar brands when situations are ambiguous due to:", "High involvement", "Cognitive shortcuts", "Brand loyalty lock-in", "Price sensitivity", "B", "Ambiguity triggers heuristic reliance.")
  ,@("In the expectancy-value model, total evaluation is computed by:", "Summing raw attribute scores", "Multiplying each belief by its importance and summing", "Averaging costs only", "Ranking alternatives randomly", "B", "Weighted sum of beliefs × importance.")
  ,@("Elimination-by-aspects heuristic works by:", "Comparing only the best alternative", "Removing options that fail on key attributes sequentially", "Calculating net present value", "Using group consensus only", "B", "Options lacking minimum attribute levels are dropped.")
  ,@("The ABC model of attitude formation stands for:", "Affect, Behavior, Cognition", "Appeal, Brand, Choice", "Awareness, Benefit, Cost", "Acquisition, Belief, Compliance", "A", "It maps feelings, actions, and thoughts.")
  ,@("Which factor does not influence postpurchase value co-creation?", "Network interactions", "In-use support", "Employee engagement", "Up-front price only", "D", "Price impacts purchase, not co-creation.")
  ,@("Physical evidence in S–O–R includes all except:", "Store layout", "Packaging design", "Employee uniforms", "Consumer motivation", "D", "Motivation is part of “Organism”.")
  ,@("Under high involvement, consumers are more likely to:", "Use heuristics exclusively", "Engage in extensive information processing", "Act on impulse only", "Ignore product attributes", "B", "High involvement demands deeper evaluation.")
  ,@("Which step of the decision process can lead to word-of-mouth or complaints?", "Evaluation of alternatives", "Purchase decision", "Postpurchase behavior", "Information search", "C", "After buying, consumers share feedback.")
  ,@("Which is not an “Example from Multimodal Advertising”?", "Scent-infused billboards", "Interactive touch-screens", "Plain text flyers", "Audio-visual kiosks", "C", "Plain text engages only vision.")
  ,@("Motivation in consumer behavior is influenced by:", "Sensation only", "Needs and goals", "Pricing solely", "Package color exclusively", "B", "It stems from unmet needs and desired ends.")
  ,@("The first consumer decision step—problem recognition—occurs when:", "A need arises from a discrepancy between actual and ideal state", "The buyer completes post-purchase evaluation", "The buyer finalizes payment method", "A value co-creation opportunity appears", "A", "Recognition of a gap triggers the process.")
  ,@("Which is an external stimulus triggering problem recognition?", "Boredom", "A friend’s recommendation", "Hunger pangs", "Personal goals", "B", "External cues come from the environment.")
)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Analysing Business and Consumer"

# Row 8 (index 7 in $data) was originally left half-filled (question text
# only) on the first pass and its answers/explanation were typed in later,
# after the rest of the table -- reproduce that same fill order so the
# shared-string table comes out identical.
$deferredRow = 7

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    if ($r -eq $deferredRow) {
        $ws.Cells.Item($r + 1, 1).Value = $row[0]
        continue
    }
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$row = $data[$deferredRow]
for ($c = 1; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($deferredRow + 1, $c + 1).Value = $row[$c]
}

# Header row formatting: reuse the same bold/centered/bordered header
# style already used on the other sheets by copying it across.
$sourceHeader = $wb.Worksheets.Item(4).Range("A1:G1")
$sourceHeader.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)

# Make the new sheet the active tab with the selection left on the row
# just below the last data row, mirroring the source file.
$ws.Activate()
$lastRow = $data.Length + 1
$ws.Range("A" + $lastRow + ":XFD" + $lastRow).Select()
